# Actualización SmartScore desde Streamlit (Laura Tamariz Valdepeña)
# Appends a new response row (row 14) to the SmartScore results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 14

# Helper: write a value as genuine text (even when it "looks" numeric,
# e.g. "0.589") without leaving a quote-prefix / number-format style
# behind on the cell, matching the plain, unstyled inline-string cells
# used throughout this sheet.
function Set-TextCell([string]$addr, [string]$text) {
    $rng = $ws.Range($addr)
    if ($text -eq "") {
        $rng.Value = ""
        return
    }
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

Set-TextCell "A$row" "Laura Tamariz Valdepeña_20251120_212917"
Set-TextCell "B$row" ""
Set-TextCell "C$row" "Laura Tamariz Valdepeña"

$ws.Range("D$row").Value = 20

Set-TextCell "E$row" "Female"
Set-TextCell "F$row" "2025-11-20 21:29:17"

$pesos = "{`n  ""portion"": 0.8,`n  ""diet"": 0.42857142857142855,`n  ""salt"": 0.4,`n  ""fat"": 0.8,`n  ""natural"": 0.6,`n  ""convenience"": 1.0,`n  ""price"": 0.6`n}"
Set-TextCell "G$row" $pesos

Set-TextCell "H$row" "Nongshim Neoguri Spicy Seafood"
Set-TextCell "I$row" "0.589"
Set-TextCell "J$row" "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"

Set-TextCell "K$row" "Maruchan Ramen Sabor Pollo"
Set-TextCell "L$row" "0.520"
Set-TextCell "M$row" "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"

Set-TextCell "N$row" "Nissin Chow Mein Teriyaki Beef"
Set-TextCell "O$row" "0.494"
Set-TextCell "P$row" "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"

Set-TextCell "Q$row" "Kraft Macaroni & Cheese Dinner"
Set-TextCell "R$row" "0.562"
Set-TextCell "S$row" "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"

Set-TextCell "T$row" "Velveeta Original Shells & Cheese (microwave cups)"
Set-TextCell "U$row" "0.553"
Set-TextCell "V$row" "Muy cremoso, porción individual, rápido, salado, ideal para niños"

Set-TextCell "W$row" "Amy’s Macaroni & Cheese (frozen)"
Set-TextCell "X$row" "0.545"
Set-TextCell "Y$row" "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"

Set-TextCell "Z$row" "Wild Planet Wild Tuna Pasta Salad"
Set-TextCell "AA$row" "0.721"
Set-TextCell "AB$row" "Sabor fresco, buena proteína, saludable, porción algo pequeña"

Set-TextCell "AC$row" "StarKist Chicken Creations (Chicken Salad)"
Set-TextCell "AD$row" "0.622"
Set-TextCell "AE$row" "Portátil, saludable, fácil, buena textura, sabor suave"

Set-TextCell "AF$row" "Jack Link’s Beef Jerky Original"
Set-TextCell "AG$row" "0.610"
Set-TextCell "AH$row" "Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña"

# The "Pesos" (G) cell holds a multi-line JSON blob. Writing multi-line
# text makes the engine mark the row with an explicit custom height;
# re-running AutoFit brings the row back to the sheet's natural/default
# height so no stray ht/customHeight attributes are left behind.
$ws.Rows.Item($row).AutoFit()
